$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "maa://24702 (94.1), maa://25390 (97.42), maa://36681 (90.77)"
$ws.Range("AA2").Value = "maa://21246 (91.2), maa://36684 (98.61), ***maa://22731 (6.67)"
$ws.Range("AE2").Value = "maa://25251 (92.5), ***maa://21730 (17.19), ***maa://39501 (23.08), *maa://36675 (60.0)"
$ws.Range("W3").Value = "maa://27396 (84.97), maa://27484 (95.74), maa://27480 (82.35)"
$ws.Range("C4").Value = "maa://24632 (93.28), **maa://24303 (36.36), maa://22499 (85.71), maa://22746 (100.0)"
$ws.Range("S6").Value = "*maa://37411 (75.0)"
$ws.Range("AE6").Value = "*maa://33152 (60.61), ***maa://22770 (28.57)"
$ws.Range("AA8").Value = "maa://25389 (88.89)"
$ws.Range("C10").Value = "***maa://25695 (19.19), **maa://32237 (37.84), ***maa://34206 (18.18), ***maa://39951 (18.52), **maa://39243 (33.33)"
$ws.Range("S11").Value = "maa://22747 (94.37), maa://22501 (98.15)"
$ws.Range("W11").Value = "maa://36713 (97.84)"
$ws.Range("C13").Value = "maa://24999 (91.46), maa://36673 (91.8), maa://25001 (85.51)"
$ws.Range("G13").Value = "*maa://21248 (75.36), **maa://22728 (47.62)"
$ws.Range("K15").Value = "*maa://21334 (52.0)"
$ws.Range("G18").Value = "maa://24421 (90.14)"
$ws.Range("K20").Value = "maa://41331 (89.19)"
$ws.Range("K23").Value = "maa://39756 (92.41), maa://39875 (95.74)"
$ws.Range("O23").Value = "maa://30587 (91.67), *maa://29748 (75.2), ***maa://29785 (15.15), *maa://37566 (78.95)"
$ws.Range("AE24").Value = "maa://22523 (85.03), *maa://36672 (76.74), maa://29910 (94.12), **maa://21440 (34.55)"
$ws.Range("AE25").Value = "maa://20108 (96.12), maa://24621 (96.4), maa://36676 (100.0), maa://22771 (84.62), maa://37772 (100.0)"
$ws.Range("W28").Value = "maa://39929 (86.29), ***maa://39723 (14.71), maa://41749 (81.25)"
$ws.Range("AE28").Value = "maa://36660 (93.8), *maa://36701 (64.0)"
$ws.Range("AE29").Value = "*maa://24080 (69.04), ***maa://34960 (9.09)"
$ws.Range("G32").Value = "maa://21895 (97.01), maa://36667 (98.11), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("S32").Value = "maa://41108 (91.18), maa://41238 (94.44)"
$ws.Range("AE34").Value = "*maa://32650 (64.29)"
$ws.Range("K35").Value = "maa://41296 (98.0)"
$ws.Range("AE38").Value = "maa://36697 (85.16)"
$ws.Range("S44").Value = "maa://39366 (83.33)"
$ws.Range("G46").Value = "maa://35931 (92.54)"
$ws.Range("G60").Value = "**maa://40438 (33.33)"
